$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "53db1b10-3d9a-4d9b-93d6-bb719553e342"
$ws.Range("B2").Value = 0.5649999999999999
$ws.Range("C2").Value = 0.131
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "Mitigation needed"

$ws.Range("A3").Value = "6f7e2563-4400-4e4b-9f39-02b0ea2de25a"
$ws.Range("B3").Value = 0.5649999999999999
$ws.Range("C3").Value = 0.131
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "Mitigation needed"
